$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Min refund calculation altered ---
# Apply a "Text" number format to the whole Min_Amount column (E2:E41).
# This mirrors the author re-typing several minimum-amount cells as text
# (so they keep leading/visual formatting) while touching the column format
# as a whole.
$ws.Range("E2:E41").NumberFormat = "@"

# Now re-enter the specific minimum-amount values that changed. Because the
# column is text-formatted, these are stored as text but Excel's arithmetic
# (F = E-1, H = E*G/100) still coerces them back to numbers for the formulas,
# while the CONCATENATE in K shows the literal text.
$ws.Range("E16").Value = "54001"
$ws.Range("E17").Value = "6000"
$ws.Range("E22").Value = "150000"
$ws.Range("E26").Value = "2000"
$ws.Range("E31").Value = "30000"
$ws.Range("E32").Value = "10000"

# --- View-state: restore the scroll position / selection left by the author ---
$ws.Activate()
try { $excel.ActiveWindow.ScrollRow = 13 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}
$ws.Range("C34").Select()
